$d = $word.ActiveDocument
$find = $d.Content.Find

# --- Merge previously-split runs back into single runs (4 places) ---
$find.Execute('El arte del engaño digital consiste en obtener información de los usuarios a través de medios como teléfonos, emails, correo tradicional o contacto directo.', $true, $false, $false, $false, $false, $true, 1, $false, 'El arte del engaño digital consiste en obtener información de los usuarios a través de medios como teléfonos, emails, correo tradicional o contacto directo.', 2) | Out-Null
$find.Execute(': se presenta cuando un supuesto representante de algún servicio pregunta por información de la cuenta del cliente.', $true, $false, $false, $false, $false, $true, 1, $false, ': se presenta cuando un supuesto representante de algún servicio pregunta por información de la cuenta del cliente.', 2) | Out-Null
$find.Execute(': consiste en colocar pendrives o memorias externas con malwares en lugares de personas escogidas que puedan infectar sus computadoras.', $true, $false, $false, $false, $false, $true, 1, $false, ': consiste en colocar pendrives o memorias externas con malwares en lugares de personas escogidas que puedan infectar sus computadoras.', 2) | Out-Null
$find.Execute(' falsos con el fin de robar información.', $true, $false, $false, $false, $false, $true, 1, $false, ' falsos con el fin de robar información.', 2) | Out-Null

# --- Append new content at the end of the document body ---
$end = $d.Content.End - 1
$rng = $d.Range($end, $end)
$rng.InsertParagraphAfter()
$p = $d.Paragraphs.Last
$p.Range.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"/>')

$end = $d.Content.End - 1
$rng = $d.Range($end, $end)
$rng.InsertParagraphAfter()
$p = $d.Paragraphs.Last
$p.Range.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"/>')

$end = $d.Content.End - 1
$rng = $d.Range($end, $end)
$rng.InsertParagraphAfter()
$p = $d.Paragraphs.Last
$p.Range.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"/>')

$end = $d.Content.End - 1
$rng = $d.Range($end, $end)
$rng.InsertParagraphAfter()
$p = $d.Paragraphs.Last
$p.Range.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"/>')

$end = $d.Content.End - 1
$rng = $d.Range($end, $end)
$rng.InsertParagraphAfter()
$p = $d.Paragraphs.Last
$p.Range.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"/>')

$end = $d.Content.End - 1
$rng = $d.Range($end, $end)
$rng.InsertParagraphAfter()
$p = $d.Paragraphs.Last
$p.Range.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="Prrafodelista"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="2"/></w:numPr></w:pPr><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:lastRenderedPageBreak/><w:t xml:space="preserve">Un software maligno también se lo conoce como </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>MALWARE</w:t></w:r><w:proofErr w:type="gramEnd"/></w:p>')

$end = $d.Content.End - 1
$rng = $d.Range($end, $end)
$rng.InsertParagraphAfter()
$p = $d.Paragraphs.Last
$p.Range.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="Prrafodelista"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="2"/></w:numPr></w:pPr><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>¿Cuál es la diferencia entre virus y troyano?</w:t></w:r></w:p>')

$end = $d.Content.End - 1
$rng = $d.Range($end, $end)
$rng.InsertParagraphAfter()
$p = $d.Paragraphs.Last
$p.Range.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="Prrafodelista"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="2"/></w:numPr></w:pPr><w:r><w:t>El virus es un componente de software que se copia a si mismo en varios lugares, mientras que el troyano es un programa sin licencia que necesita de la ejecución del usuario</w:t></w:r></w:p>')

$end = $d.Content.End - 1
$rng = $d.Range($end, $end)
$rng.InsertParagraphAfter()
$p = $d.Paragraphs.Last
$p.Range.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="Prrafodelista"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="2"/></w:numPr></w:pPr><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t xml:space="preserve">Los </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>s</w:t></w:r><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>pywares</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t xml:space="preserve"> NO dañ</w:t></w:r><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>an los d</w:t></w:r><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>ispositivos</w:t></w:r></w:p>')

$end = $d.Content.End - 1
$rng = $d.Range($end, $end)
$rng.InsertParagraphAfter()
$p = $d.Paragraphs.Last
$p.Range.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="Prrafodelista"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="2"/></w:numPr></w:pPr><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t xml:space="preserve">¿En qué se diferencian los </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>rootkits</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t xml:space="preserve"> de las demás amenazas?</w:t></w:r></w:p>')

$end = $d.Content.End - 1
$rng = $d.Range($end, $end)
$rng.InsertParagraphAfter()
$p = $d.Paragraphs.Last
$p.Range.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="Prrafodelista"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="2"/></w:numPr></w:pPr><w:r><w:t xml:space="preserve">Los </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>rootkits</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> logran esconderse de los softwares antimalware o antivirus</w:t></w:r></w:p>')

$end = $d.Content.End - 1
$rng = $d.Range($end, $end)
$rng.InsertParagraphAfter()
$p = $d.Paragraphs.Last
$p.Range.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="Prrafodelista"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="2"/></w:numPr></w:pPr><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>¿Cuáles son los aspectos de la información que se deben proteger para evitar ataques?</w:t></w:r></w:p>')

$end = $d.Content.End - 1
$rng = $d.Range($end, $end)
$rng.InsertParagraphAfter()
$p = $d.Paragraphs.Last
$p.Range.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="Prrafodelista"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="2"/></w:numPr></w:pPr><w:r><w:t>Confidencialidad, integridad y disponibilidad</w:t></w:r></w:p>')
